# Updated cryptos list on Mon Oct 28 09:42:37 UTC 2024 with GitHub Actions
#
# Refreshes the per-coin Price (column D) and Volume(1h) (column E) values
# on the single worksheet, and fixes the row order for ARBITRUM / BabyDogeCoin
# (rows 48-49 swap places along with their Link/Price/Volume data).
#
# Note: several Price values are digit-only strings (e.g. "591.08", "0.517",
# "1.00") that must stay plain text, exactly like the source cells
# (t="inlineStr"), instead of being auto-coerced to numbers by Excel. We use
# the standard Excel "leading apostrophe" text-entry marker for those so the
# stored/displayed value is still the plain digit string, not a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "68.254.57"
$ws.Range("E2").Value = "  +1.69%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "2.507.35"
$ws.Range("E3").Value = "  +1.61%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.05%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "'591.08"
$ws.Range("E5").Value = "  +1.48%  "

# --- Row 6: Solana ---
$ws.Range("E6").Value = "  +1.40%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  -0.01%  "

# --- Row 8: XRP ---
$ws.Range("D8").Value = "'0.517"

# --- Row 9: LidoStakedEther ---
$ws.Range("D9").Value = "2.507.47"
$ws.Range("E9").Value = "  +1.60%  "

# --- Row 10: Dogecoin ---
$ws.Range("D10").Value = "'0.143"
$ws.Range("E10").Value = "  +3.95%  "

# --- Row 11: TRON ---
$ws.Range("E11").Value = "  -1.00%  "

# --- Row 12: Toncoin ---
$ws.Range("E12").Value = "  +0.70%  "

# --- Row 13: Cardano ---
$ws.Range("D13").Value = "'0.336"
$ws.Range("E13").Value = "  +0.97%  "

# --- Row 14: WrappedliquidstakedEther2.0 ---
$ws.Range("D14").Value = "2.991.39"

# --- Row 15: Avalanche ---
$ws.Range("D15").Value = "'25.81"
$ws.Range("E15").Value = "  +1.76%  "

# --- Row 16: WrappedBTC ---
$ws.Range("D16").Value = "68.120.35"
$ws.Range("E16").Value = "  +1.69%  "

# --- Row 17: ShibaInu ---
$ws.Range("E17").Value = "  +0.37%  "

# --- Row 18: WrappedEther ---
$ws.Range("D18").Value = "2.499.08"
$ws.Range("E18").Value = "  +1.47%  "

# --- Row 19: Chainlink ---
$ws.Range("D19").Value = "'10.98"

# --- Row 20: Uniswap ---
$ws.Range("D20").Value = "'7.43"
$ws.Range("E20").Value = "  -0.27%  "

# --- Row 21: BitcoinCash ---
$ws.Range("D21").Value = "'351.01"
$ws.Range("E21").Value = "  +0.62%  "

# --- Row 22: Polkadot ---
$ws.Range("E22").Value = "  +4.94%  "

# --- Row 23: Litecoin ---
$ws.Range("D23").Value = "'71.25"
$ws.Range("E23").Value = "  +2.86%  "

# --- Row 24: Dai ---
$ws.Range("E24").Value = "  -0.11%  "

# --- Row 25: NEARProtocol ---
$ws.Range("E25").Value = "  +0.28%  "

# --- Row 26: SuiNetwork ---
$ws.Range("E26").Value = "  -3.99%  "

# --- Row 27: Aptos ---
$ws.Range("D27").Value = "'9.11"
$ws.Range("E27").Value = "  -0.04%  "

# --- Row 28: WrappedeETH ---
$ws.Range("E28").Value = "  +1.82%  "

# --- Row 29: Binance-PegBSC-USD ---
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.09%  "

# --- Row 30: PEPE ---
$ws.Range("E30").Value = "  -0.43%  "

# --- Row 31: Bittensor ---
$ws.Range("D31").Value = "'510.14"
$ws.Range("E31").Value = "  +2.09%  "

# --- Row 32: InternetComputer(DFINITY) ---
$ws.Range("D32").Value = "'7.77"
$ws.Range("E32").Value = "  +0.68%  "

# --- Row 33: Fetch.AI ---
$ws.Range("D33").Value = "'1.25"
$ws.Range("E33").Value = "  +2.10%  "

# --- Row 34: PancakeSwap ---
$ws.Range("E34").Value = "  +1.03%  "

# --- Row 35: FirstDigitalUSD ---
$ws.Range("E35").Value = "  -0.02%  "

# --- Row 36: Kaspa ---
$ws.Range("E36").Value = "  +1.09%  "

# --- Row 37: Monero ---
$ws.Range("D37").Value = "'162.05"
$ws.Range("E37").Value = "  +0.02%  "

# --- Row 38: WhiteBITCoin ---
$ws.Range("D38").Value = "'18.68"
$ws.Range("E38").Value = "  +0.04%  "

# --- Row 39: EthereumClassic ---
$ws.Range("D39").Value = "'18.33"
$ws.Range("E39").Value = "  +1.14%  "

# --- Row 40: ImmutableX ---
$ws.Range("E40").Value = "  -0.12%  "

# --- Row 42: Stacks ---
$ws.Range("E42").Value = "  +3.29%  "

# --- Row 43: PolygonEcosystemToken ---
$ws.Range("D43").Value = "'0.327"
$ws.Range("E43").Value = "  +0.30%  "

# --- Row 44: RenderToken ---
$ws.Range("E44").Value = "  -0.18%  "

# --- Row 45: dogwifhat ---
$ws.Range("E45").Value = "  +1.96%  "

# --- Row 46: Aave ---
$ws.Range("D46").Value = "'150.87"
$ws.Range("E46").Value = "  +6.24%  "

# --- Row 47: Filecoin ---
$ws.Range("E47").Value = "  +2.60%  "

# --- Row 48: was BabyDogeCoin -> now ARBITRUM ---
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.518"
$ws.Range("E48").Value = "  +1.46%  "

# --- Row 49: was ARBITRUM -> now BabyDogeCoin ---
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0258"
$ws.Range("E49").Value = "  +1.92%  "

# --- Row 50: Optimism ---
$ws.Range("E50").Value = "  +1.73%  "
